$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 35
$ws.Range("B3").Value = 35
$ws.Range("B4").Value = 9
$ws.Range("B5").Value = 9
$ws.Range("B6").Value = 6
$ws.Range("B7").Value = 6

$ws.Range("B8").Select()
